$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C16").Value = "1047446850"
$ws.Range("D16").Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Range("E16").Value = "1909"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 1300000

$ws.Range("C17").Value = "1047446850"
$ws.Range("D17").Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Range("E17").Value = "1908"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 1300000

$ws.Range("C18").Value = "1047446850"
$ws.Range("D18").Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Range("E18").Value = "1907"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 1300000

$ws.Range("C19").Value = "1047446850"
$ws.Range("D19").Value = "MIRLLAN YULIETH MARTINEZ PITALUA"
$ws.Range("E19").Value = "1906"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 1300000

$ws.Range("C20").Value = "45442352"
$ws.Range("D20").Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Range("E20").Value = "1908"
$ws.Range("F20").Value = 8833
$ws.Range("G20").Value = 828116

$ws.Range("C21").Value = "45442352"
$ws.Range("D21").Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Range("E21").Value = "1907"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 828116

$ws.Range("C22").Value = "45442352"
$ws.Range("D22").Value = "GLADYS DEL SOCORRO LOPEZ LLERENA"
$ws.Range("E22").Value = "1906"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116
